# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across all 8 sheets per the scheduled price-refresh diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 24807256
$ws.Range("I70").Value = 17864688
$ws.Range("J70").Value = 31749822
$ws.Range("K70").Value = 53594064
$ws.Range("L70").Value = 95249466
$ws.Range("M70").Value = -53593794
$ws.Range("N70").Value = -95250006
$ws.Range("H73").Value = 24807256
$ws.Range("I73").Value = 17864688
$ws.Range("J73").Value = 31749822
$ws.Range("K73").Value = 53594064
$ws.Range("L73").Value = 95249466
$ws.Range("M73").Value = -53593128
$ws.Range("N73").Value = -95251338
$ws.Range("H98").Value = 55559436
$ws.Range("I98").Value = 55559436
$ws.Range("K98").Value = 55559436
$ws.Range("M98").Value = -55557938
$ws.Range("H100").Value = 1146.6666
$ws.Range("I100").Value = 1152.0952
$ws.Range("J100").Value = 1134
$ws.Range("K100").Value = 1152.0952
$ws.Range("L100").Value = 1134
$ws.Range("M100").Value = -611.0952
$ws.Range("N100").Value = -2216
$ws.Range("H106").Value = 33334616
$ws.Range("J106").Value = 925
$ws.Range("L106").Value = 925
$ws.Range("N106").Value = -2187
$ws.Range("H122").Value = 55559436
$ws.Range("I122").Value = 55559436
$ws.Range("K122").Value = 166678308
$ws.Range("M122").Value = -166675858
$ws.Range("H131").Value = 1165.8334
$ws.Range("I131").Value = 1165.8334
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3497.5002
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 1542.4998
$ws.Range("H133").Value = 8441922
$ws.Range("J133").Value = 8441922
$ws.Range("L133").Value = 8441922
$ws.Range("N133").Value = -8452042
$ws.Range("H137").Value = 5708.654
$ws.Range("I137").Value = 3919.4
$ws.Range("J137").Value = 8148.5454
$ws.Range("K137").Value = 11758.2
$ws.Range("L137").Value = 24445.6362
$ws.Range("M137").Value = -9208.200000000001
$ws.Range("N137").Value = -29545.6362
$ws.Range("N131").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1956869.1
$ws.Range("I32").Value = 1956869.1
$ws.Range("K32").Value = 1956869.1
$ws.Range("M32").Value = -1956582.1
$ws.Range("H61").Value = 71437350
$ws.Range("I61").Value = 1275.8572
$ws.Range("J61").Value = 142873420
$ws.Range("K61").Value = 1275.8572
$ws.Range("L61").Value = 142873420
$ws.Range("M61").Value = -1063.8572
$ws.Range("N61").Value = -142873844
$ws.Range("H97").Value = 16671631
$ws.Range("I97").Value = 6270
$ws.Range("K97").Value = 6270
$ws.Range("M97").Value = -5774
$ws.Range("H102").Value = 7694355.5
$ws.Range("I102").Value = 10528039
$ws.Range("J102").Value = 2928.4285
$ws.Range("K102").Value = 10528039
$ws.Range("L102").Value = 2928.4285
$ws.Range("M102").Value = -10526417
$ws.Range("N102").Value = -6172.4285
$ws.Range("H112").Value = 40763.2
$ws.Range("J112").Value = 40763.2
$ws.Range("L112").Value = 40763.2
$ws.Range("N112").Value = -43717.2
$ws.Range("H122").Value = 15574.1
$ws.Range("I122").Value = 26577.111
$ws.Range("J122").Value = 6571.636
$ws.Range("K122").Value = 79731.333
$ws.Range("L122").Value = 19714.908
$ws.Range("M122").Value = -77281.333
$ws.Range("N122").Value = -24614.908
$ws.Range("H132").Value = 7152.528
$ws.Range("I132").Value = 4279.1055
$ws.Range("J132").Value = 10364
$ws.Range("K132").Value = 12837.3165
$ws.Range("L132").Value = 31092
$ws.Range("M132").Value = -10307.3165
$ws.Range("N132").Value = -36152
$ws.Range("H136").Value = 71437350
$ws.Range("I136").Value = 1275.8572
$ws.Range("J136").Value = 142873420
$ws.Range("K136").Value = 3827.5716
$ws.Range("L136").Value = 428620260
$ws.Range("M136").Value = -1277.5716
$ws.Range("N136").Value = -428625360

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 26789902
$ws.Range("J107").Value = 5433
$ws.Range("L107").Value = 5433
$ws.Range("N107").Value = -9273
$ws.Range("H134").Value = 5004523.5
$ws.Range("I134").Value = 7814026
$ws.Range("K134").Value = 23442078
$ws.Range("M134").Value = -23439543

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2400
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H22").Value = 959.5625
$ws.Range("I22").Value = 866.75
$ws.Range("J22").Value = 1238
$ws.Range("K22").Value = 866.75
$ws.Range("L22").Value = 1238
$ws.Range("M22").Value = -516.75
$ws.Range("N22").Value = -1938
$ws.Range("H31").Value = 5966.3213
$ws.Range("I31").Value = 2346.4348
$ws.Range("J31").Value = 8489.272000000001
$ws.Range("K31").Value = 2346.4348
$ws.Range("L31").Value = 8489.272000000001
$ws.Range("M31").Value = -2051.4348
$ws.Range("N31").Value = -9079.272000000001
$ws.Range("H34").Value = 5966.3213
$ws.Range("I34").Value = 2346.4348
$ws.Range("J34").Value = 8489.272000000001
$ws.Range("K34").Value = 2346.4348
$ws.Range("L34").Value = 8489.272000000001
$ws.Range("M34").Value = -2144.4348
$ws.Range("N34").Value = -8893.272000000001
$ws.Range("H111").Value = 91296.664
$ws.Range("J111").Value = 91296.664
$ws.Range("L111").Value = 91296.664
$ws.Range("N111").Value = -99476.664
$ws.Range("H122").Value = 38069.895
$ws.Range("I122").Value = 1986.5
$ws.Range("K122").Value = 5959.5
$ws.Range("M122").Value = -3509.5
$ws.Range("H141").Value = 167777.25
$ws.Range("J141").Value = 167777.25
$ws.Range("L141").Value = 167777.25
$ws.Range("N141").Value = -178137.25
$ws.Range("N4").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4242551.5
$ws.Range("J4").Value = 964164.5600000001
$ws.Range("L4").Value = 2892493.68
$ws.Range("N4").Value = -2892717.68
$ws.Range("H55").Value = 5563380.5
$ws.Range("J55").Value = 6258554
$ws.Range("L55").Value = 18775662
$ws.Range("N55").Value = -18776016
$ws.Range("H113").Value = 2584.1
$ws.Range("J113").Value = 3384.4285
$ws.Range("L113").Value = 10153.2855
$ws.Range("N113").Value = -14493.2855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 186856.98
$ws.Range("I113").Value = 529186.75
$ws.Range("J113").Value = 6182.9443
$ws.Range("K113").Value = 529186.75
$ws.Range("L113").Value = 6182.9443
$ws.Range("M113").Value = -527016.75
$ws.Range("N113").Value = -10522.9443
$ws.Range("H122").Value = 3625638.5
$ws.Range("I122").Value = 4832613.5
$ws.Range("J122").Value = 4713.6
$ws.Range("K122").Value = 14497840.5
$ws.Range("L122").Value = 14140.8
$ws.Range("M122").Value = -14495390.5
$ws.Range("N122").Value = -19040.8
$ws.Range("H132").Value = 6565.2
$ws.Range("J132").Value = 11657
$ws.Range("L132").Value = 34971
$ws.Range("N132").Value = -40031

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 60000
$ws.Range("J2").Value = 62500
$ws.Range("L2").Value = 62500
$ws.Range("N2").Value = -62724
$ws.Range("H16").Value = 1499.75
$ws.Range("I16").Value = 1499.75
$ws.Range("K16").Value = 1499.75
$ws.Range("M16").Value = -1329.75
$ws.Range("H22").Value = 10709.167
$ws.Range("I22").Value = 808.3333
$ws.Range("K22").Value = 808.3333
$ws.Range("M22").Value = -513.3333
$ws.Range("H27").Value = 10709.167
$ws.Range("I27").Value = 808.3333
$ws.Range("K27").Value = 808.3333
$ws.Range("M27").Value = -701.3333
$ws.Range("H136").Value = 9161.4
$ws.Range("I136").Value = 3276.3333
$ws.Range("J136").Value = 13575.2
$ws.Range("K136").Value = 9828.999899999999
$ws.Range("L136").Value = 40725.60000000001
$ws.Range("M136").Value = -7278.999899999999
$ws.Range("N136").Value = -45825.60000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 49721.668
$ws.Range("I2").Value = 49721.668
$ws.Range("K2").Value = 49721.668
$ws.Range("M2").Value = -49609.668
$ws.Range("H41").Value = 16305.714
$ws.Range("J41").Value = 16305.714
$ws.Range("L41").Value = 16305.714
$ws.Range("N41").Value = -17085.714
$ws.Range("H96").Value = 2895.75
$ws.Range("J96").Value = 4794
$ws.Range("L96").Value = 4794
$ws.Range("N96").Value = -7540
$ws.Range("H107").Value = 12346600
$ws.Range("I107").Value = 317.2
$ws.Range("J107").Value = 27779454
$ws.Range("K107").Value = 951.5999999999999
$ws.Range("L107").Value = 83338362
$ws.Range("M107").Value = 968.4000000000001
$ws.Range("N107").Value = -83342202
$ws.Range("H136").Value = 25900228
$ws.Range("I136").Value = 34483830
$ws.Range("K136").Value = 103451490
$ws.Range("M136").Value = -103448940
